$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking values that must remain literal text
# (matching the source inlineStr cells), so force Text format before assigning.
$deCells = @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","E48","D49","E49","E50","E51")
foreach ($addr in $deCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$values = @{
    'D2' = '330.34'
    'E2' = '2.90%'
    'D3' = '41.08'
    'E3' = '4.55%'
    'D4' = '5.683'
    'E4' = '-4.01%'
    'D5' = '0.08172'
    'E5' = '1.96%'
    'E6' = '7.24%'
    'D7' = '8.756'
    'E7' = '1.29%'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '2.952'
    'E8' = '0.39%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9252'
    'E9' = '-0.93%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1250'
    'E10' = '-0.19%'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1959'
    'E11' = '0.20%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09295'
    'E12' = '1.59%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03696'
    'E13' = '5.17%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.1054'
    'E14' = '10.11%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001299'
    'E15' = '0.47%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.006150'
    'E16' = '-0.79%'
    'B17' = 'HotbitToken'
    'C17' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D17' = '0.004290'
    'E17' = '-1.75%'
    'D18' = '3.419'
    'E18' = '1.83%'
    'B19' = 'GateToken'
    'C19' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D19' = '4.545'
    'E19' = '-0.46%'
    'B20' = 'BitpandaEcosystemToken'
    'C20' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D20' = '0.3485'
    'E20' = '-1.43%'
    'B21' = 'MCDex'
    'C21' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D21' = '8.338'
    'E21' = '-5.20%'
    'B22' = 'ProBitToken'
    'C22' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D22' = '0.1381'
    'E22' = '-3.39%'
    'B23' = 'ZBToken'
    'C23' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D23' = '0.2654'
    'E23' = '10.00%'
    'B24' = 'CoinExToken'
    'C24' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D24' = '0.04449'
    'E24' = '-0.02%'
    'B25' = 'BitKan'
    'C25' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D25' = '0.001270'
    'E25' = '0.65%'
    'E26' = '3.65%'
    'D39' = '0.02782'
    'E39' = '15.87%'
    'D40' = '0.05484'
    'E40' = '5.84%'
    'D41' = '0.007679'
    'E41' = '3.31%'
    'D42' = '0.009439'
    'E42' = '2.28%'
    'D43' = '0.1415'
    'E43' = '0.94%'
    'D44' = '0.002124'
    'E44' = '0.13%'
    'D45' = '0.01212'
    'E45' = '8.54%'
    'D46' = '0.00006857'
    'E46' = '1.90%'
    'E47' = '0.13%'
    'E48' = '60.53%'
    'D49' = '0.003532'
    'E49' = '17.46%'
    'E50' = '0.13%'
    'E51' = '0.13%'
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
